$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row 4 (dimension grows from A1:H3 to A1:H4)
$ws.Range("A4").Value = "####1er Torneo Federativo - C.A.E. - Sub 23, Prejuveniles y sub 23 (28 de Febrero y 1 de Marzo) - Juniors (Domingo 1 de Marzo)"
$ws.Range("B4").Value = "Juveniles"
$ws.Range("C4").Value = "caballeros"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "Liberatori, Augusto"
$ws.Range("F4").Value = 83
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = 83
